$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for Price/Volume columns so numeric-looking
# strings (e.g. "20.237.33") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '20.237.33', '  +1.17%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.451.28', '  +2.90%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.016', '  +1.62%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '277.53', '  +1.30%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.8999', '  -9.97%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.3683', '  -0.52%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.3132', '  +2.02%  '),
    @(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '39.10', '  -0.69%  '),
    @(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '1.019', '  +1.87%  '),
    @(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06461', '  -1.63%  '),
    @(12, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.010', '  +0.94%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.387', '  -0.51%  '),
    @(14, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '17.38', '  +2.36%  '),
    @(15, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.471.92', '  +4.68%  '),
    @(16, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.087', '  -1.53%  '),
    @(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001013', '  +0.60%  '),
    @(18, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.05607', '  -2.84%  '),
    @(19, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9048', '  -9.49%  '),
    @(20, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '67.26', '  -8.75%  '),
    @(21, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.469', '  -2.56%  '),
    @(22, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '14.31', '  -1.11%  '),
    @(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '11.06', '  +1.90%  '),
    @(24, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '2.274', '  -1.70%  '),
    @(25, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '20.442.48', '  +2.21%  '),
    @(26, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '2.188', '  -3.80%  '),
    @(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '135.59', '  -2.21%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '16.93', '  +0.35%  '),
    @(29, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.633.66', '  +4.35%  '),
    @(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '109.84', '  +0.72%  '),
    @(31, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '3.812', '  -0.34%  '),
    @(32, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.8010', '  -6.30%  '),
    @(33, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.857', '  -9.71%  '),
    @(34, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.07700', '  -0.02%  '),
    @(35, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05941', '  +2.45%  '),
    @(36, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '1.443', '  +12.67%  '),
    @(37, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.727', '  -1.73%  '),
    @(38, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '1.146', '  +7.69%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.02002', '  -2.42%  '),
    @(40, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '10.21', '  -0.94%  '),
    @(41, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1822', '  -5.41%  '),
    @(42, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '0.9172', '  -8.32%  '),
    @(43, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.559', '  +0.72%  '),
    @(44, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.5255', '  -0.83%  '),
    @(45, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '12.07', '  -0.47%  '),
    @(46, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.663', '  -21.10%  '),
    @(47, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '120.16', '  +9.06%  '),
    @(48, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.5118', '  -0.15%  '),
    @(49, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '1.760', '  -2.62%  '),
    @(50, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.06341', '  +2.73%  '),
    @(51, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9972', '  -0.31%  ')
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).Value = $item[4]
}

$wb.Save()
